# Update the "Date" value on the Metadata sheet and the "System URI" values
# on the three "Include" sheets, per the commit's re-organisation of the
# published documentation site (ansforge.github.io -> interop.esante.gouv.fr).

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update the generation Date string (B8) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-07-21T11:52:46+00:00"

# --- Include #0 sheet: TRE-R249-Sexe system URI (B4) ---
$inc0 = $wb.Worksheets.Item("Include #0")
$inc0.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R249-Sexe"

# --- Include #1 sheet: TRE-R267-SexeProvenanceISO system URI (B4) ---
$inc1 = $wb.Worksheets.Item("Include #1")
$inc1.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R267-SexeProvenanceISO"

# --- Include #2 sheet: TRE-R303-HL7v3AdministrativeGender system URI (B4) ---
$inc2 = $wb.Worksheets.Item("Include #2")
$inc2.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R303-HL7v3AdministrativeGender"
